$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" summary text (Binance rate lines) ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$text = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 15.15 = 62575.76 pesos`n✅ 62575.76 pesos = 15.16 = 979.55 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $text

# --- tasas: update the N10/O10/N12/O12 rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 66
$ws2.Range("O10").Value = 4130
$ws2.Range("N12").Value = 4126.78
$ws2.Range("O12").Value = 64.59999999999999
